$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - column F ("想去人数" / interested-count) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 8100
$ws1.Range("F4").Value = 2710
$ws1.Range("F7").Value = 808
$ws1.Range("F8").Value = 624
$ws1.Range("F9").Value = 104
$ws1.Range("F10").Value = 70
$ws1.Range("F12").Value = 874
$ws1.Range("F13").Value = 3322
$ws1.Range("F14").Value = 222
$ws1.Range("F15").Value = 111
$ws1.Range("F16").Value = 757
$ws1.Range("F21").Value = 283
$ws1.Range("F23").Value = 355
$ws1.Range("F25").Value = 135
$ws1.Range("F27").Value = 296
$ws1.Range("F28").Value = 34
$ws1.Range("F32").Value = 585
$ws1.Range("F35").Value = 21
$ws1.Range("F38").Value = 109

# Sheet "全部类型" (All Types) - same column F updates (mirrors 展览 rows, offset by 2)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 8100
$ws4.Range("F6").Value = 2710
$ws4.Range("F9").Value = 808
$ws4.Range("F10").Value = 624
$ws4.Range("F11").Value = 104
$ws4.Range("F12").Value = 70
$ws4.Range("F14").Value = 874
$ws4.Range("F16").Value = 3322
$ws4.Range("F17").Value = 222
$ws4.Range("F18").Value = 111
$ws4.Range("F20").Value = 757
$ws4.Range("F26").Value = 283
$ws4.Range("F28").Value = 355
$ws4.Range("F30").Value = 135
$ws4.Range("F32").Value = 296
$ws4.Range("F33").Value = 34
$ws4.Range("F37").Value = 585
$ws4.Range("F40").Value = 21
$ws4.Range("F43").Value = 109
